$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 207, shifting existing rows 207..255 down to 208..256.
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new record.
$ws.Cells.Item(207, 1).Value2  = 4
$ws.Cells.Item(207, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(207, 3).Value2  = "Los Lagos"
$ws.Cells.Item(207, 4).Value2  = 44641
$ws.Cells.Item(207, 5).Value2  = 10
$ws.Cells.Item(207, 6).Value2  = 100112040
$ws.Cells.Item(207, 7).Value2  = "Cilantro"
$ws.Cells.Item(207, 8).Value2  = "Sin especificar"
$ws.Cells.Item(207, 9).Value2  = "Primera"
$ws.Cells.Item(207, 10).Value2 = 70
$ws.Cells.Item(207, 11).Value2 = 5000
$ws.Cells.Item(207, 12).Value2 = 5000
$ws.Cells.Item(207, 13).Value2 = 5000
$ws.Cells.Item(207, 14).Value2 = '$/docena de atados (2 kilos)'
$ws.Cells.Item(207, 15).Value2 = 'Región de La Araucanía'
$ws.Cells.Item(207, 16).Value2 = 2500
$ws.Cells.Item(207, 17).Value2 = 2
$ws.Cells.Item(207, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same date number-format style as the
# other date cells in column D (style index 2 on the sheet).
$ws.Cells.Item(207, 4).NumberFormat = $ws.Cells.Item(208, 4).NumberFormat
